{"js": "// The {{r single_to_double_newlines(value(clientmessage.text)) | paragraphs | markdown }}\n// merge field is being switched to the simpler {{p value(clientmessage.text) }} field\n// used by the new docx templates (see commit message \"switched to new docx templates\").\n//\n// We perform the edit as three small, targeted search-and-replace operations so that\n// only the text that actually changed is touched; everything else in the paragraph\n// (including the `_GoBack` bookmark) is left exactly where it was.\n\nconst body = context.document.body;\n\n// 1) \"{{r \" -> \"{{p \"\nlet startMarker = body.search(\"{{r \", { matchCase: true, matchWholeWord: false });\nstartMarker.load(\"items\");\nawait context.sync();\nif (startMarker.items.length > 0) {\n  startMarker.items[0].insertText(\"{{p \", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) drop the \"single_to_double_newlines(\" wrapper call\nlet wrapperOpen = body.search(\"single_to_double_newlines(\", { matchCase: true, matchWholeWord: false });\nwrapperOpen.load(\"items\");\nawait context.sync();\nif (wrapperOpen.items.length > 0) {\n  wrapperOpen.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 3) drop the matching \")\" plus the \" | paragraphs | markdown\" filter pipeline\nlet wrapperTail = body.search(\") | paragraphs | markdown\", { matchCase: true, matchWholeWord: false });\nwrapperTail.load(\"items\");\nawait context.sync();\nif (wrapperTail.items.length > 0) {\n  wrapperTail.items[0].insertText(\"\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# The {{r single_to_double_newlines(value(clientmessage.text)) | paragraphs | markdown }}\n# merge field is being switched to the simpler {{p value(clientmessage.text) }} field used\n# by the new docx templates (see commit message \"switched to new docx templates\").\n#\n# We perform the edit as a handful of small, targeted Find/Replace passes over\n# $d.Content so that only the text that actually changed is touched; everything\n# else in the paragraph (including the `_GoBack` bookmark) stays right where it was.\n\n$d = $word.ActiveDocument\n\n# 1) \"{{r \" -> \"{{p \"\n$d.Content.Find.Execute(\"{{r \", $false, $false, $false, $false, $false, $true, 1, $false, \"{{p \", 2) | Out-Null\n\n# 2) drop the \"single_to_double_newlines(\" wrapper call\n$d.Content.Find.Execute(\"single_to_double_newlines(\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# 3) drop \") | paragraphs | \" (everything up to, but not past, the \"_GoBack\" bookmark)\n$d.Content.Find.Execute(\") | paragraphs | \", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n\n# 4) drop the trailing \"markdown\" filter name (the part after the bookmark)\n$d.Content.Find.Execute(\"markdown\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 2) | Out-Null\n"}
